# Insert the new "Our Tested hypothisis" slide as slide 13 (1-based index),
# i.e. right after the current slide 12 (sldId 312) and before the old
# slide 13 (sldId 308). Layout 6 = "Title and long list" (ctrTitle + body),
# the same layout already used by e.g. slide 2 ("Questions we explored").
$p = $ppt.ActivePresentation
$s = $p.Slides.Add(13, 6)

# --- Title placeholder ------------------------------------------------
$title = $s.Shapes.Item(1)
$title.Name = "Title 1"
$titleTr = $title.TextFrame.TextRange
$titleTr.Text = "Our Tested "
$titleTr.LanguageID = "en-US"
$titleRun2 = $titleTr.InsertAfter("hypothisis")
$titleRun2.LanguageID = "en-US"

# --- Body placeholder ---------------------------------------------------
$body = $s.Shapes.Item(2)
$body.Name = "Text Placeholder 2"
$bodyTr = $body.TextFrame.TextRange

$lines = @(
  "Hypothesis Testing:",
  "",
  "Hypothesis: Highly dense continents are mostly affected by Covid-19 due to closer proximity of potential carriers. ",
  "",
  "",
  "Analysis:",
  "",
  "Our calculated Degree of Freedom = 5",
  "Confidence Interval = 95%",
  "Critical Value = 11.070497693516351",
  "TMP_BLANK_BULLET",
  "Conclusion:",
  "",
  "The Chi-square value exceeds the critical value, which means the statistically significant but we cannot accept the Null Hypothesis because the P value is 0",
  "Therefore, Population density of the continent does not contribute to the spread of "
)
$bodyTr.Text = [string]::Join([char]13, $lines)
$bodyTr.LanguageID = "en-US"

# Paragraphs that read as plain section headers (no bullet), indented
# like the rest of the deck's "label:" lines.
$noBulletParas = @(1, 2, 6, 7, 12, 13)
foreach ($idx in $noBulletParas) {
    $para = $bodyTr.Paragraphs($idx, 1)
    $para.ParagraphFormat.Bullet.Visible = 0
}

# Paragraphs that are regular bullet list items.
$bulletParas = @(3, 8, 9, 10, 11, 14, 15)
foreach ($idx in $bulletParas) {
    $para = $bodyTr.Paragraphs($idx, 1)
    $pf = $para.ParagraphFormat
    $pf.Bullet.Font.Name = "Arial"
    $pf.Bullet.Visible = 1
    $pf.Bullet.Type = 1
    $pf.Bullet.Character = 8226
}

# Paragraph 11 is a blank bullet line (formatting kept, text cleared).
$blankBullet = $bodyTr.Paragraphs(11, 1)
$blankBullet.Text = ""

# Paragraph 15 ends with a separate "Covid" run (own rPr, no dirty flag
# carried over from proofing, same as the rest of the deck's split runs).
$lastPara = $bodyTr.Paragraphs(15, 1)
$covidRun = $lastPara.InsertAfter("Covid")
$covidRun.LanguageID = "en-US"
